$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update project names / approval status per "Updates to project approvals"
$ws.Range("D10").Value = "Double pendulum"

$ws.Range("D8").Value = "Library of Lighting Models"

$ws.Range("D9").Value = "Optimal Solar Panel Tilt"
$ws.Range("G9").Value = "Yes"

$ws.Range("G10").Value = "Yes"

# Move active cell selection to D10 as in the final saved state
$ws.Range("D10").Select()
